# Generate Report for Handback
# Adds a second handback record (63878f5d-2025-4b80-aa58-1777943db14a) alongside
# the existing one (now renamed from b607217c... to 612d4ba2...), across the
# Overview / zh-cn / de-de sheets, and bumps the "Latest HO Xliff Generate Date".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper colour (BGR int) matching the workbook's existing HyperLink font
# (styles.xml fonts: <u/><color rgb="FF6495ED"/>)
# ---------------------------------------------------------------------------
$hyperlinkColor = 15570276

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# -- existing row 2: the handback file was renamed (b607217c... -> 612d4ba2...)
$ws1.Range("A2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md"

$ws1.Range("B2").Value = "e2e\612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md"
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a3e6c4956a045711a893c9bb570daaf62b27c7a/e2e/612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md", "", "", "e2e\612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md") | Out-Null
$ws1.Range("B2").Font.Underline = $true
$ws1.Range("B2").Font.Color = $hyperlinkColor

# -- existing row 2: bump the generate-date text --------------------------
$ws1.Range("G2").Value = "2016-09-04 19:06:31"
$ws1.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- new row 3: second handback file ---------------------------------------
$ws1.Range("A3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.md"

$ws1.Range("B3").Value = "e2e\63878f5d-2025-4b80-aa58-1777943db14a.md"
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a3e6c4956a045711a893c9bb570daaf62b27c7a/e2e/63878f5d-2025-4b80-aa58-1777943db14a.md", "", "", "e2e\63878f5d-2025-4b80-aa58-1777943db14a.md") | Out-Null
$ws1.Range("B3").Font.Underline = $true
$ws1.Range("B3").Font.Color = $hyperlinkColor

$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-09-04 19:06:31"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$lo1 = $ws1.ListObjects.Item(1)
$lo1.Resize($ws1.Range("A1:G3"))

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")

# -- existing row 2: rename source/target/handoff/handback file names ------
$ws2.Range("A2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a3e6c4956a045711a893c9bb570daaf62b27c7a/e2e/612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md", "", "", "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md") | Out-Null
$ws2.Range("A2").Font.Underline = $true
$ws2.Range("A2").Font.Color = $hyperlinkColor

$ws2.Range("G2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.beda524222e96ebd3d4f6d891550b2ca8bdd86ae.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-04 19:06:27"
$ws2.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("I2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md"
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/71954e4985e0b8c5c138dc571730df479a89fdd8/e2e/612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md", "", "", "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md") | Out-Null
$ws2.Range("I2").Font.Underline = $true
$ws2.Range("I2").Font.Color = $hyperlinkColor

$ws2.Range("J2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.beda524222e96ebd3d4f6d891550b2ca8bdd86ae.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-04 19:06:45"
$ws2.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- new row 3: second handback file ----------------------------------------
$ws3row = $ws2.Range("A3:P3")

$ws2.Range("A3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.md"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a3e6c4956a045711a893c9bb570daaf62b27c7a/e2e/63878f5d-2025-4b80-aa58-1777943db14a.md", "", "", "63878f5d-2025-4b80-aa58-1777943db14a.md") | Out-Null
$ws2.Range("A3").Font.Underline = $true
$ws2.Range("A3").Font.Color = $hyperlinkColor

$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.a4caaa0ff1cddf8b0a6c40e76bfd6652e9982c00.zh-cn.xlf"

$ws2.Range("H3").Value = "2016-09-04 19:06:27"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("I3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.md"
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/71954e4985e0b8c5c138dc571730df479a89fdd8/e2e/63878f5d-2025-4b80-aa58-1777943db14a.md", "", "", "63878f5d-2025-4b80-aa58-1777943db14a.md") | Out-Null
$ws2.Range("I3").Font.Underline = $true
$ws2.Range("I3").Font.Color = $hyperlinkColor

$ws2.Range("J3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.a4caaa0ff1cddf8b0a6c40e76bfd6652e9982c00.zh-cn.xlf"

$ws2.Range("K3").Value = "2016-09-04 19:06:45"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

$lo2 = $ws2.ListObjects.Item(1)
$lo2.Resize($ws2.Range("A1:P3"))

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")

# -- existing row 2: rename source/target/handoff/handback file names ------
$ws3.Range("A2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a3e6c4956a045711a893c9bb570daaf62b27c7a/e2e/612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md", "", "", "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md") | Out-Null
$ws3.Range("A2").Font.Underline = $true
$ws3.Range("A2").Font.Color = $hyperlinkColor

$ws3.Range("G2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.beda524222e96ebd3d4f6d891550b2ca8bdd86ae.de-de.xlf"

$ws3.Range("H2").Value = "2016-09-04 19:06:31"
$ws3.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("I2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md"
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cc1242acc70b13231915944ad27ab80d4fe6e660/e2e/612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md", "", "", "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.md") | Out-Null
$ws3.Range("I2").Font.Underline = $true
$ws3.Range("I2").Font.Color = $hyperlinkColor

$ws3.Range("J2").Value = "612d4ba2-bc7d-4b9a-a1d4-8e97ba5cf1ba.beda524222e96ebd3d4f6d891550b2ca8bdd86ae.de-de.xlf"

$ws3.Range("K2").Value = "2016-09-04 19:06:53"
$ws3.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# -- new row 3: second handback file ----------------------------------------
$ws3.Range("A3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.md"
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9a3e6c4956a045711a893c9bb570daaf62b27c7a/e2e/63878f5d-2025-4b80-aa58-1777943db14a.md", "", "", "63878f5d-2025-4b80-aa58-1777943db14a.md") | Out-Null
$ws3.Range("A3").Font.Underline = $true
$ws3.Range("A3").Font.Color = $hyperlinkColor

$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.a4caaa0ff1cddf8b0a6c40e76bfd6652e9982c00.de-de.xlf"

$ws3.Range("H3").Value = "2016-09-04 19:06:31"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("I3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.md"
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cc1242acc70b13231915944ad27ab80d4fe6e660/e2e/63878f5d-2025-4b80-aa58-1777943db14a.md", "", "", "63878f5d-2025-4b80-aa58-1777943db14a.md") | Out-Null
$ws3.Range("I3").Font.Underline = $true
$ws3.Range("I3").Font.Color = $hyperlinkColor

$ws3.Range("J3").Value = "63878f5d-2025-4b80-aa58-1777943db14a.a4caaa0ff1cddf8b0a6c40e76bfd6652e9982c00.de-de.xlf"

$ws3.Range("K3").Value = "2016-09-04 19:06:53"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

$lo3 = $ws3.ListObjects.Item(1)
$lo3.Resize($ws3.Range("A1:P3"))
